# Edit script for NORTH_CAROLINA_2022.xlsx
# - Rename header row columns to snake_case field names
# - Title-case municipality/state name strings (e.g. "de" -> "De", "del" -> "Del")
# - Fix 12 D-column percentage cells that differ by 1 floating-point ULP
# - Remove trailing footer/notes rows (1524-1529)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (A1:D1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Title-case municipality / state names in columns A and B ---
    $ws.Range("B7").Value = 'Pabellón De Arteaga'
    $ws.Range("B8").Value = 'Rincón De Romos'
    $ws.Range("B9").Value = 'San Francisco De Los Romo'
    $ws.Range("B33").Value = 'Amatenango De La Frontera'
    $ws.Range("B36").Value = 'Bejucal De Ocampo'
    $ws.Range("B38").Value = 'Benemérito De Las Américas'
    $ws.Range("B45").Value = 'Chiapa De Corzo'
    $ws.Range("B51").Value = 'Comitán De Domínguez'
    $ws.Range("B77").Value = 'Marqués De Comillas'
    $ws.Range("B78").Value = 'Mazapa De Madero'
    $ws.Range("B84").Value = 'Ocozocoautla De Espinosa'
    $ws.Range("B97").Value = 'Salto De Agua'
    $ws.Range("B98").Value = 'San Cristóbal De Las Casas'
    $ws.Range("B101").Value = 'Santiago El Pinar'
    $ws.Range("B134").Value = 'Hidalgo Del Parral'
    $ws.Range("B138").Value = 'San Francisco Del Oro'
    $ws.Range("B156").Value = 'San Juan De Sabinas'
    $ws.Range("A168").Value = 'Ciudad De México'
    $ws.Range("B172").Value = 'Cuajimalpa De Morelos'
    $ws.Range("B195").Value = 'Nombre De Dios'
    $ws.Range("B198").Value = 'Pánuco De Coronado'
    $ws.Range("B203").Value = 'San Juan De Guadalupe'
    $ws.Range("B204").Value = 'San Juan Del Río'
    $ws.Range("A211").Value = 'Estado De México'
    $ws.Range("B211").Value = 'Acambay De Ruíz Castañeda'
    $ws.Range("B214").Value = 'Almoloya De Alquisiras'
    $ws.Range("B215").Value = 'Almoloya De Juárez'
    $ws.Range("B221").Value = 'Atizapán De Zaragoza'
    $ws.Range("B232").Value = 'Coacalco De Berriozábal'
    $ws.Range("B238").Value = 'Ecatepec De Morelos'
    $ws.Range("B244").Value = 'Ixtapan De La Sal'
    $ws.Range("B257").Value = 'Naucalpan De Juárez'
    $ws.Range("B265").Value = 'San Felipe Del Progreso'
    $ws.Range("B267").Value = 'San Simón De Guerrero'
    $ws.Range("B277").Value = 'Tenango Del Valle'
    $ws.Range("B285").Value = 'Tlalnepantla De Baz'
    $ws.Range("B290").Value = 'Valle De Bravo'
    $ws.Range("B291").Value = 'Valle De Chalco Solidaridad'
    $ws.Range("B292").Value = 'Villa De Allende'
    $ws.Range("B293").Value = 'Villa Del Carbón'
    $ws.Range("B305").Value = 'Apaseo El Alto'
    $ws.Range("B306").Value = 'Apaseo El Grande'
    $ws.Range("B314").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
    $ws.Range("B318").Value = 'Jaral Del Progreso'
    $ws.Range("B326").Value = 'Purísima Del Rincón'
    $ws.Range("B330").Value = 'San Diego De La Unión'
    $ws.Range("B332").Value = 'San Francisco Del Rincón'
    $ws.Range("B334").Value = 'San Luis De La Paz'
    $ws.Range("B335").Value = 'Santa Cruz De Juventino Rosas'
    $ws.Range("B336").Value = 'Silao De La Victoria'
    $ws.Range("B341").Value = 'Valle De Santiago'
    $ws.Range("B347").Value = 'Acapulco De Juárez'
    $ws.Range("B350").Value = 'Ajuchitlán Del Progreso'
    $ws.Range("B351").Value = 'Alcozauca De Guerrero'
    $ws.Range("B355").Value = 'Atenango Del Río'
    $ws.Range("B357").Value = 'Atoyac De Álvarez'
    $ws.Range("B358").Value = 'Ayutla De Los Libres'
    $ws.Range("B361").Value = 'Buenavista De Cuéllar'
    $ws.Range("B362").Value = 'Chilapa De Álvarez'
    $ws.Range("B363").Value = 'Chilpancingo De Los Bravo'
    $ws.Range("B364").Value = 'Coahuayutla De José María Izazaga'
    $ws.Range("B369").Value = 'Coyuca De Benítez'
    $ws.Range("B370").Value = 'Coyuca De Catalán'
    $ws.Range("B374").Value = 'Cuetzala Del Progreso'
    $ws.Range("B375").Value = 'Cutzamala De Pinzón'
    $ws.Range("B381").Value = 'Huitzuco De Los Figueroa'
    $ws.Range("B382").Value = 'Iguala De La Independencia'
    $ws.Range("B384").Value = 'Ixcateopan De Cuauhtémoc'
    $ws.Range("B385").Value = 'Zihuatanejo De Azueta'
    $ws.Range("B387").Value = 'La Unión De Isidoro Montes De Oca'
    $ws.Range("B390").Value = 'Mártir De Cuilapan'
    $ws.Range("B403").Value = 'Taxco De Alarcón'
    $ws.Range("B405").Value = 'Técpan De Galeana'
    $ws.Range("B407").Value = 'Tepecoacuilco De Trujano'
    $ws.Range("B409").Value = 'Tixtla De Guerrero'
    $ws.Range("B412").Value = 'Tlalixtaquilla De Maldonado'
    $ws.Range("B413").Value = 'Tlapa De Comonfort'
    $ws.Range("B425").Value = 'Agua Blanca De Iturbide'
    $ws.Range("B431").Value = 'Atotonilco De Tula'
    $ws.Range("B432").Value = 'Atotonilco El Grande'
    $ws.Range("B437").Value = 'Cuautepec De Hinojosa'
    $ws.Range("B440").Value = 'Huasca De Ocampo'
    $ws.Range("B443").Value = 'Huejutla De Reyes'
    $ws.Range("B446").Value = 'Jacala De Ledezma'
    $ws.Range("B452").Value = 'Mineral Del Chico'
    $ws.Range("B453").Value = 'Mineral Del Monte'
    $ws.Range("B454").Value = 'Mixquiahuala De Juárez'
    $ws.Range("B455").Value = 'Molango De Escamilla'
    $ws.Range("B457").Value = 'Nopala De Villagrán'
    $ws.Range("B458").Value = 'Omitlán De Juárez'
    $ws.Range("B459").Value = 'Pachuca De Soto'
    $ws.Range("B462").Value = 'Progreso De Obregón'
    $ws.Range("B468").Value = 'Santiago De Anaya'
    $ws.Range("B469").Value = 'Santiago Tulantepec De Lugo Guerrero'
    $ws.Range("B473").Value = 'Tenango De Doria'
    $ws.Range("B475").Value = 'Tepehuacán De Guerrero'
    $ws.Range("B476").Value = 'Tepeji Del Río De Ocampo'
    $ws.Range("B479").Value = 'Tezontepec De Aldama'
    $ws.Range("B487").Value = 'Tula De Allende'
    $ws.Range("B488").Value = 'Tulancingo De Bravo'
    $ws.Range("B489").Value = 'Villa De Tezontepec'
    $ws.Range("B493").Value = 'Zacualtipán De Ángeles'
    $ws.Range("B500").Value = 'Atotonilco El Alto'
    $ws.Range("B502").Value = 'Autlán De Navarro'
    $ws.Range("B515").Value = 'Encarnación De Díaz'
    $ws.Range("B520").Value = 'Huejuquilla El Alto'
    $ws.Range("B521").Value = 'Ixtlahuacán Del Río'
    $ws.Range("B524").Value = 'Jilotlán De Los Dolores'
    $ws.Range("B527").Value = 'Lagos De Moreno'
    $ws.Range("B532").Value = 'Ojuelos De Jalisco'
    $ws.Range("B537").Value = 'San Juan De Los Lagos'
    $ws.Range("B538").Value = 'San Martín De Bolaños'
    $ws.Range("B539").Value = 'San Sebastián Del Oeste'
    $ws.Range("B540").Value = 'Santa María De Los Ángeles'
    $ws.Range("B542").Value = 'Tamazula De Gordiano'
    $ws.Range("B544").Value = 'Techaluta De Montenegro'
    $ws.Range("B547").Value = 'Tepatitlán De Morelos'
    $ws.Range("B548").Value = 'Tizapán El Alto'
    $ws.Range("B549").Value = 'Tlajomulco De Zúñiga'
    $ws.Range("B553").Value = 'Unión De San Antonio'
    $ws.Range("B554").Value = 'Valle De Guadalupe'
    $ws.Range("B559").Value = 'Yahualica De González Gallo'
    $ws.Range("B560").Value = 'Zacoalco De Torres'
    $ws.Range("B563").Value = 'Zapotitlán De Vadillo'
    $ws.Range("B564").Value = 'Zapotlán El Grande'
    $ws.Range("B585").Value = 'Coalcomán De Vázquez Pallares'
    $ws.Range("B644").Value = 'Tiquicheo De Nicolás Romero'
    $ws.Range("B665").Value = 'Coatlán Del Río'
    $ws.Range("B676").Value = 'Puente De Ixtla'
    $ws.Range("B682").Value = 'Tlaltizapán De Zapata'
    $ws.Range("B690").Value = 'Bahía De Banderas'
    $ws.Range("B694").Value = 'Ixtlán Del Río'
    $ws.Range("B701").Value = 'Santa María Del Oro'
    $ws.Range("B722").Value = 'Mier Y Noriega'
    $ws.Range("B723").Value = 'Montemorelos'
    $ws.Range("B725").Value = 'San Nicolás De Los Garza'
    $ws.Range("B729").Value = 'Acatlán De Pérez Figueroa'
    $ws.Range("B736").Value = 'Chiquihuitlán De Benito Juárez'
    $ws.Range("B737").Value = 'Ciénega De Zimatlán'
    $ws.Range("B741").Value = 'Constancia Del Rosario'
    $ws.Range("B743").Value = 'Cuilápam De Guerrero'
    $ws.Range("B744").Value = 'Fresnillo De Trujano'
    $ws.Range("B745").Value = 'Guadalupe De Ramírez'
    $ws.Range("B746").Value = 'Guelatao De Juárez'
    $ws.Range("B747").Value = 'Guevea De Humboldt'
    $ws.Range("B748").Value = 'Heroica Ciudad De Ejutla De Crespo'
    $ws.Range("B749").Value = 'Heroica Ciudad De Huajuapan De León'
    $ws.Range("B750").Value = 'Heroica Ciudad De Tlaxiaco'
    $ws.Range("B751").Value = 'Huautla De Jiménez'
    $ws.Range("B752").Value = 'Ixtlán De Juárez'
    $ws.Range("B753").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
    $ws.Range("B760").Value = 'Mariscala De Juárez'
    $ws.Range("B761").Value = 'Mártires De Tacubaya'
    $ws.Range("B763").Value = 'Mazatlán Villa De Flores'
    $ws.Range("B765").Value = 'Miahuatlán De Porfirio Díaz'
    $ws.Range("B766").Value = 'Mixistlán De La Reforma'
    $ws.Range("B768").Value = 'Nejapa De Madero'
    $ws.Range("B770").Value = 'Oaxaca De Juárez'
    $ws.Range("B771").Value = 'Ocotlán De Morelos'
    $ws.Range("B772").Value = 'Pinotepa De Don Luis'
    $ws.Range("B773").Value = 'Putla Villa De Guerrero'
    $ws.Range("B774").Value = 'Reforma De Pineda'
    $ws.Range("B778").Value = 'San Agustín De Las Juntas'
    $ws.Range("B784").Value = 'San Antonino El Alto'
    $ws.Range("B786").Value = 'San Antonio De La Cal'
    $ws.Range("B797").Value = 'San Felipe Jalapa De Díaz'
    $ws.Range("B803").Value = 'San Francisco Del Mar'
    $ws.Range("B820").Value = 'San Juan Bautista Lo De Soto'
    $ws.Range("B830").Value = 'San Juan Del Estado'
    $ws.Range("B831").Value = 'San Juan Del Río'
    $ws.Range("B861").Value = 'San Miguel Del Puerto'
    $ws.Range("B862").Value = 'San Miguel Del Río'
    $ws.Range("B863").Value = 'San Miguel El Grande'
    $ws.Range("B876").Value = 'San Pablo Villa De Mitla'
    $ws.Range("B880").Value = 'San Pedro El Alto'
    $ws.Range("B911").Value = 'Santa Cruz De Bravo'
    $ws.Range("B969").Value = 'Santo Domingo De Morelos'
    $ws.Range("B983").Value = 'Tamazulápam Del Espíritu Santo'
    $ws.Range("B984").Value = 'Tataltepec De Valdés'
    $ws.Range("B985").Value = 'Teotitlán De Flores Magón'
    $ws.Range("B986").Value = 'Tezoatlán De Segura Y Luna'
    $ws.Range("B987").Value = 'Tlacolula De Matamoros'
    $ws.Range("B988").Value = 'Totontepec Villa De Morelos'
    $ws.Range("B990").Value = 'Villa De Etla'
    $ws.Range("B991").Value = 'Villa De Tututepec De Melchor Ocampo'
    $ws.Range("B992").Value = 'Villa Sola De Vega'
    $ws.Range("B994").Value = 'Zapotitlán Del Río'
    $ws.Range("B997").Value = 'Zimatlán De Álvarez'
    $ws.Range("B1015").Value = 'Ayotoxco De Guerrero'
    $ws.Range("B1018").Value = 'Chalchicomula De Sesma'
    $ws.Range("B1035").Value = 'Cuayuca De Andrade'
    $ws.Range("B1036").Value = 'Cuetzalan Del Progreso'
    $ws.Range("B1048").Value = 'Huehuetlán El Chico'
    $ws.Range("B1052").Value = 'Ixcamilpa De Guerrero'
    $ws.Range("B1054").Value = 'Izúcar De Matamoros'
    $ws.Range("B1063").Value = 'Los Reyes De Juárez'
    $ws.Range("B1070").Value = 'Palmar De Bravo'
    $ws.Range("B1092").Value = 'Tecali De Herrera'
    $ws.Range("B1097").Value = 'Tepanco De López'
    $ws.Range("B1098").Value = 'Tepatlaxco De Hidalgo'
    $ws.Range("B1104").Value = 'Tepexi De Rodríguez'
    $ws.Range("B1106").Value = 'Tetela De Ocampo'
    $ws.Range("B1110").Value = 'Tlacotepec De Benito Juárez'
    $ws.Range("B1122").Value = 'Tuzamapan De Galeana'
    $ws.Range("B1138").Value = 'Amealco De Bonfil'
    $ws.Range("B1140").Value = 'Cadereyta De Montes'
    $ws.Range("B1146").Value = 'Jalpan De Serra'
    $ws.Range("B1147").Value = 'Landa De Matamoros'
    $ws.Range("B1150").Value = 'Pinal De Amoles'
    $ws.Range("B1153").Value = 'San Juan Del Río'
    $ws.Range("B1164").Value = 'Axtla De Terrazas'
    $ws.Range("B1169").Value = 'Cerro De San Pedro'
    $ws.Range("B1170").Value = 'Ciudad Del Maíz'
    $ws.Range("B1179").Value = 'Mexquitic De Carmona'
    $ws.Range("B1184").Value = 'San Ciro De Acosta'
    $ws.Range("B1188").Value = 'Santa María Del Río'
    $ws.Range("B1190").Value = 'Soledad De Graciano Sánchez'
    $ws.Range("B1200").Value = 'Villa De Arista'
    $ws.Range("B1201").Value = 'Villa De Arriaga'
    $ws.Range("B1202").Value = 'Villa De Guadalupe'
    $ws.Range("B1203").Value = 'Villa De La Paz'
    $ws.Range("B1204").Value = 'Villa De Ramos'
    $ws.Range("B1205").Value = 'Villa De Reyes'
    $ws.Range("B1231").Value = 'Nacozari De García'
    $ws.Range("B1244").Value = 'Jalpa De Méndez'
    $ws.Range("B1279").Value = 'Soto La Marina'
    $ws.Range("B1291").Value = 'Contla De Juan Cuamatzi'
    $ws.Range("B1295").Value = 'Ixtacuixtla De Mariano Matamoros'
    $ws.Range("B1299").Value = 'San Pablo Del Monte'
    $ws.Range("B1315").Value = 'Alto Lucero De Gutiérrez Barrios'
    $ws.Range("B1319").Value = 'Amatlán De Los Reyes'
    $ws.Range("B1329").Value = 'Boca Del Río'
    $ws.Range("B1331").Value = 'Camarón De Tejeda'
    $ws.Range("B1334").Value = 'Castillo De Teayo'
    $ws.Range("B1342").Value = 'Chinampa De Gorostiza'
    $ws.Range("B1355").Value = 'Cosamaloapan De Carpio'
    $ws.Range("B1370").Value = 'Hueyapan De Ocampo'
    $ws.Range("B1371").Value = 'Huiloapan De Cuauhtémoc'
    $ws.Range("B1372").Value = 'Ignacio De La Llave'
    $ws.Range("B1376").Value = 'Ixhuatlán De Madero'
    $ws.Range("B1377").Value = 'Ixhuatlán Del Café'
    $ws.Range("B1378").Value = 'Ixhuatlán Del Sureste'
    $ws.Range("B1388").Value = 'Juchique De Ferrer'
    $ws.Range("B1390").Value = 'Landero Y Coss'
    $ws.Range("B1393").Value = 'Lerdo De Tejada'
    $ws.Range("B1397").Value = 'Martínez De La Torre'
    $ws.Range("B1400").Value = 'Medellín De Bravo'
    $ws.Range("B1404").Value = 'Nanchital De Lázaro Cárdenas Del Río'
    $ws.Range("B1414").Value = 'Paso De Ovejas'
    $ws.Range("B1415").Value = 'Paso Del Macho'
    $ws.Range("B1418").Value = 'Poza Rica De Hidalgo'
    $ws.Range("B1426").Value = 'Sayula De Alemán'
    $ws.Range("B1429").Value = 'Soledad De Doblado'
    $ws.Range("B1436").Value = 'Tatahuicapan De Juárez'
    $ws.Range("B1453").Value = 'Tlacotepec De Mejía'
    $ws.Range("B1466").Value = 'Vega De Alatorre'
    $ws.Range("B1489").Value = 'Concepción Del Oro'
    $ws.Range("B1497").Value = 'Jiménez Del Teul'
    $ws.Range("B1502").Value = 'Mezquital Del Oro'
    $ws.Range("B1505").Value = 'Nochistlán De Mejía'
    $ws.Range("B1506").Value = 'Noria De Ángeles'
    $ws.Range("B1515").Value = 'Teúl De González Ortega'
    $ws.Range("B1516").Value = 'Tlaltenango De Sánchez Román'
    $ws.Range("B1518").Value = 'Villa De Cos'

# --- 3. Fix floating point 1-ULP differences in column D (recomputed percentages) ---
$ws.Range("D73").Value = 0.0009406296685663696
$ws.Range("D167").Value = 0.0009406296685663696
$ws.Range("D302").Value = 0.0009406296685663696
$ws.Range("D325").Value = 0.0009406296685663696
$ws.Range("D375").Value = 0.009018978586842252
$ws.Range("D622").Value = 0.0009406296685663696
$ws.Range("D715").Value = 0.0009406296685663696
$ws.Range("D849").Value = 0.0009406296685663696
$ws.Range("D952").Value = 0.0009406296685663696
$ws.Range("D964").Value = 0.0009406296685663696
$ws.Range("D1088").Value = 0.0009406296685663696
$ws.Range("D1192").Value = 0.0009406296685663696
$ws.Range("D1485").Value = 0.0009406296685663696

# --- 4. Remove trailing footer/notes rows 1524-1529 ---
$ws.Range("A1524:A1529").EntireRow.Delete()
